$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Lab Exercise")
$ws2 = $wb.Worksheets.Item("Assignment")

# --- Fix the AIC formulas: the correction term numerator should be the
# sample size (N) for each block (20, 40, 60, 80) instead of a hardcoded 1 ---

# Block 1 (N = 20): B16 single formula, B17:B24 shared formula
$ws1.Range("B16").Formula = "=-2*LN(B2) + (2*1*(20/(20-1-1)))"
$ws1.Range("B17:B24").Formula = "=-2*LN(B3) + (2*1*(20/(20-1-1)))"

# Block 2 (N = 40): B28 single formula, B29:B36 shared formula
$ws1.Range("B28").Formula = "=-2*LN(C2) + (2*1*(40/(40-1-1)))"
$ws1.Range("B29:B36").Formula = "=-2*LN(C3) + (2*1*(40/(40-1-1)))"

# Block 3 (N = 60): B40 single formula, B41:B48 shared formula
$ws1.Range("B40").Formula = "=-2*LN(D2) + (2*1*(60/(60-1-1)))"
$ws1.Range("B41:B48").Formula = "=-2*LN(D3) + (2*1*(60/(60-1-1)))"

# Block 4 (N = 80): B53 single formula, B54:B61 shared formula
$ws1.Range("B53").Formula = "=-2*LN(E2) + (2*1*(80/(80-1-1)))"
$ws1.Range("B54:B61").Formula = "=-2*LN(E3) + (2*1*(80/(80-1-1)))"

# --- Update active sheet / selection state ---
# "Lab Exercise" becomes the active (selected) tab with selection at G55
# "Assignment" is no longer the selected tab, selection stays at G2
$ws1.Activate()
$ws1.Range("G55").Select()

$ws2.Range("G2").Select()

$ws1.Activate()
